# "break out stock.yaml completed"
#
# 1) On the "day" sheet, columns D864:D875 were stored as inline-string text
#    ("500495", etc.) and must become true numeric values (500495, ...).
# 2) On the "week" sheet, 29 new data rows (530..557 in 1-based Excel terms,
#    i.e. rows 529-557 in the underlying XML) are appended after the
#    existing last row (528), extending the used range to A1:I557.

$wb = $excel.ActiveWorkbook

# --- 1) Fix bsecode column type on "day" sheet (rows 864-875) ---------------
$day = $wb.Worksheets.Item("day")

$bsecodes = @{
    864 = 500495
    865 = 500331
    866 = 500696
    867 = 500790
    868 = 500257
    869 = 533309
    870 = 524494
    871 = 532424
    872 = 500302
    873 = 540133
    874 = 540222
    875 = 500875
}

foreach ($r in $bsecodes.Keys) {
    $day.Cells.Item($r, 4).Value = $bsecodes[$r]
}

# --- 2) Append new rows to "week" sheet --------------------------------------
$week = $wb.Worksheets.Item("week")

$newRows = @(
    @(1,  "SHREECEM",   "Shree Cements Limited",                            500387, -1.02,               24525.05, 17563,    "week", "08/11/2024 11:34:20"),
    @(2,  "BAJFINANCE",  "Bajaj Finance Limited",                           500034, -0.12,               6895.95,  424841,   "week", "08/11/2024 11:34:20"),
    @(3,  "BRITANNIA",   "Britannia Industries Limited",                    500825, 1.02,                5747.15,  306281,   "week", "08/11/2024 11:34:20"),
    @(4,  "NAVINFLUOR",  "Navin Fluorine International Limited",            532504, 0.07000000000000001, 3543.35,  204063,   "week", "08/11/2024 11:34:20"),
    @(5,  "PIDILITIND",  "Pidilite Industries Limited",                     500331, -1.03,               3125.7,   107369,   "week", "08/11/2024 11:34:20"),
    @(6,  "COLPAL",      "Colgate Palmolive (india) Limited",               500830, -1.41,               2870.9,   482004,   "week", "08/11/2024 11:34:20"),
    @(7,  "GRASIM",      "Grasim Industries Limited",                       500300, -0.7,                2544.85,  558298,   "week", "08/11/2024 11:34:20"),
    @(8,  "ACC",         "Acc Limited",                                     500410, -1.26,               2291.4,   120875,   "week", "08/11/2024 11:34:20"),
    @(9,  "SUNPHARMA",   "Sun Pharmaceuticals Industries Limited",          524715, 1.04,                1809.6,   1644075,  "week", "08/11/2024 11:34:20"),
    @(10, "TATACOMM",    "Tata Communications Limited",                    500483, -2.2,                1766,     129546,   "week", "08/11/2024 11:34:20"),
    @(11, "TECHM",       "Tech Mahindra Limited",                           532755, 1.84,                1681.35,  1883257,  "week", "08/11/2024 11:34:20"),
    @(12, "JSWSTEEL",    "Jsw Steel Limited",                               500228, 0.67,                996.4,    3568198,  "week", "08/11/2024 11:34:20"),
    @(13, "ICICIPRULI",  "Icici Prudential Life Insurance Company Limited",540133, -0.49,               710.35,   489065,   "week", "08/11/2024 11:34:20"),
    @(14, "HINDALCO",    "Hindalco Industries Limited",                     500440, 0.36,                650.45,   7886228,  "week", "08/11/2024 11:34:20"),
    @(15, "MARICO",      "Marico Limited",                                  531642, -0.29,               629.85,   457415,   "week", "08/11/2024 11:34:20"),
    @(16, "AMBUJACEM",   "Ambuja Cements Limited",                          500425, -1.02,               564.9,    901270,   "week", "08/11/2024 11:34:20"),
    @(17, "UPL",         "Upl Limited",                                     512070, -1.68,               557.6,    696687,   "week", "08/11/2024 11:34:20"),
    @(18, "CHAMBLFERT",  "Chambal Fertilizers & Chemicals Limited",         500085, -6.4,                482.65,   4368887,  "week", "08/11/2024 11:34:20"),
    @(19, "VEDL",        "Vedanta Limited",                                 500295, 0,                   457.9,    13143777, "week", "08/11/2024 11:34:20"),
    @(20, "NTPC",        "Ntpc Limited",                                    532555, -1.52,               397.65,   7536509,  "week", "08/11/2024 11:34:20"),
    @(21, "BPCL",        "Bharat Petroleum Corporation Limited",            500547, -1.44,               310.45,   9304376,  "week", "08/11/2024 11:34:20"),
    @(22, "HINDCOPPER",  "Hindustan Copper Limited",                        513599, -3.52,               280.95,   4299350,  "week", "08/11/2024 11:34:20"),
    @(23, "NMDC",        "Nmdc Limited",                                    526371, -1.23,               235.97,   23074306, "week", "08/11/2024 11:34:20"),
    @(24, "GAIL",        "Gail (india) Limited",                            532155, -2.97,               204.17,   13570458, "week", "08/11/2024 11:34:20"),
    @(25, "ABCAPITAL",   "Aditya Birla Capital Ltd",                        540691, -1.06,               201.11,   3688004,  "week", "08/11/2024 11:34:20"),
    @(26, "TATASTEEL",   "Tata Steel Limited",                              500470, -2.24,               147.57,   61980989,"week", "08/11/2024 11:34:20"),
    @(27, "LTF",         "L&T Finance Ltd",                                 533519, -2.8,                140.11,   7010267,  "week", "08/11/2024 11:34:20"),
    @(28, "CANBK",       "Canara Bank",                                     532483, -1.31,               103.69,   19035406, "week", "08/11/2024 11:34:20"),
    @(29, "IDFCFIRSTB",  "IDFC First Bank Ltd",                             539437, -1.34,               65.63,    24563138, "week", "08/11/2024 11:34:20")
)

$startRow = 529

$r = $startRow
foreach ($row in $newRows) {
    $week.Cells.Item($r, 1).Value = $row[0]
    $week.Cells.Item($r, 2).Value = $row[1]
    $week.Cells.Item($r, 3).Value = $row[2]
    $week.Cells.Item($r, 4).Value = $row[3]
    $week.Cells.Item($r, 5).Value = $row[4]
    $week.Cells.Item($r, 6).Value = $row[5]
    $week.Cells.Item($r, 7).Value = $row[6]
    $week.Cells.Item($r, 8).Value = $row[7]
    $week.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}
